$wb = $excel.ActiveWorkbook

# The "Metadata" worksheet holds Property/Value pairs in columns A/B.
$ws = $wb.Worksheets.Item("Metadata")

# Update the "Date" property value (row 8) to the new timestamp.
$ws.Range("B8").Value = "2025-10-02T18:31:12+01:00"

# Set the "Case Sensitive" property value (row 20) to "true".
# Force text formatting so the literal word "true" is stored as a string,
# not auto-converted to an Excel boolean.
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "true"
